$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bordered/centered data-row formatting (currently only on rows 2-5)
# down through row 12 before writing any new values.
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C12").PasteSpecial(-4122)

# Update the title/date header in A1 (Ballgorithm/ESPN headers in B1/C1 stay the same)
$ws.Range("A1").Value = "NBA, Tuesday 27th Feb 2024"

# New matchup data (column A), Ballgorithm predictions (column B)
# and ESPN predictions (column C) for rows 2-12.
$aVals = @(
    "Dallas Mavericks (33-24) vs Cleveland Cavaliers (37-19)",
    "Brooklyn Nets (22-35) vs Orlando Magic (32-26)",
    "Golden State Warriors (29-27) vs Washington Wizards (9-48)",
    "Utah Jazz (27-31) vs Atlanta Hawks (25-32)",
    "Philadelphia 76ers (33-24) vs Boston Celtics (45-12)",
    "San Antonio Spurs (11-47) vs Minnesota Timberwolves (40-17)",
    "New Orleans Pelicans (34-24) vs New York Knicks (35-23)",
    "Detroit Pistons (8-49) vs Chicago Bulls (27-30)",
    "Charlotte Hornets (15-42) vs Milwaukee Bucks (37-21)",
    "Houston Rockets (25-32) vs Oklahoma City Thunder (40-17)",
    "Miami Heat (32-25) vs Portland Trail Blazers (15-41)"
)
$bVals = @(
    "Cleveland Cavaliers (65.52%)",
    "Orlando Magic (69.23%)",
    "Golden State Warriors (51.61%)",
    "Utah Jazz (62.07%)",
    "Boston Celtics (89.66%)",
    "Minnesota Timberwolves (76.92%)",
    "New York Knicks (68.97%)",
    "Chicago Bulls (53.57%)",
    "Milwaukee Bucks (76.67%)",
    "Oklahoma City Thunder (79.31%)",
    "Miami Heat (53.57%)"
)
$cVals = @(
    "Cleveland Cavaliers (71.9%)",
    "Orlando Magic (77.6%)",
    "Golden State Warriors (82.5%)",
    "Utah Jazz (52.0%)",
    "Boston Celtics (90.0%)",
    "Minnesota Timberwolves (90.3%)",
    "New York Knicks (50.8%)",
    "Chicago Bulls (77.0%)",
    "Milwaukee Bucks (90.1%)",
    "Oklahoma City Thunder (83.1%)",
    "Miami Heat (72.7%)"
)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $aVals[$i]
}
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $bVals[$i]
}
for ($i = 0; $i -lt $cVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $cVals[$i]
}

$ws.Range("A1").Select()
